# NBA_Index.xlsx — "Completed NBA_model data set and saved as NBA_model.csv"
#
# The Orlando Magic row (row 23) had its Team Name cell ("Orlando Magic")
# cleared/overwritten so it now just reads the abbreviation "ORL", matching
# column A. Once that text is gone, "Orlando Magic" is no longer referenced
# anywhere in the sheet, so Excel drops it from the shared-strings table on
# save (handled automatically by the engine).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B23").Value = "ORL"

# Reflect the author's final on-screen view state: zoomed to 130% with
# cell E16 selected when the workbook was saved.
$ws.Range("E16").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130
